$wb = $excel.ActiveWorkbook

# The workbook contains 4 sheets: 展览, 演出, 本地生活, 全部类型
# Only 展览 and 全部类型 contain data; both need identical updates to
# column F (想去人数) for several rows.

$updates = @{
    3  = 11295
    4  = 10589
    7  = 754
    12 = 10510
    18 = 94
    19 = 405
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
